$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Update role values: "PI" -> "principal investigator", "Technician" -> "technician"
$ws.Range("G2").Value = "principal investigator"
$ws.Range("G4").Value = "technician"
